# weekly_breakdown_summary.xlsx refresh
# Dashboard optimizations: roll the weekly breakdown table forward by one week
# (new "Jul 26 - Aug 01" week added, oldest "Jun 28 - Jul 04" week dropped) and
# refresh the recomputed Tickers/Winners/Losers/Win_Rate/Avg_Return/Best/Worst
# figures for every remaining week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Period (column B): shift each week's label down a row, newest on top ---
$ws.Range("B2").Value = "Jul 26 - Aug 01"
$ws.Range("B3").Value = "Jul 19 - Jul 25"
$ws.Range("B4").Value = "Jul 12 - Jul 18"
$ws.Range("B5").Value = "Jul 05 - Jul 11"

# --- Tickers (column C) ---
$ws.Range("C2").Value = 46
$ws.Range("C3").Value = 60
$ws.Range("C4").Value = 101
$ws.Range("C5").Value = 121

# --- Winners (column D) ---
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 9
$ws.Range("D4").Value = 21
$ws.Range("D5").Value = 41

# --- Losers (column E) ---
$ws.Range("E2").Value = 36
$ws.Range("E3").Value = 51
$ws.Range("E4").Value = 80
$ws.Range("E5").Value = 80

# --- Win_Rate (column F) ---
$ws.Range("F2").Value = 21.73913043478261
$ws.Range("F3").Value = 15
$ws.Range("F4").Value = 20.79207920792079
$ws.Range("F5").Value = 33.88429752066116

# --- Avg_Return (column G) ---
$ws.Range("G2").Value = -2.504204763949144
$ws.Range("G3").Value = -5.003652024233992
$ws.Range("G4").Value = -4.727910665635076
$ws.Range("G5").Value = -3.694411749866992

# --- Best_Ticker (column H) ---
$ws.Range("H2").Value = "JSL"
$ws.Range("H3").Value = "SHYAMMETL"
$ws.Range("H4").Value = "AGI"
$ws.Range("H5").Value = "ANANDRATHI"

# --- Best_Return (column I) ---
$ws.Range("I2").Value = 2.045093221563806
$ws.Range("I3").Value = 4.643413769011655
$ws.Range("I4").Value = 12.16232227488151
$ws.Range("I5").Value = 18.49389910763066

# --- Worst_Ticker (column J) ---
$ws.Range("J2").Value = "AARTIIND"
$ws.Range("J3").Value = "MRPL"
$ws.Range("J4").Value = "LODHA"
$ws.Range("J5").Value = "MSUMI"

# --- Worst_Return (column K) ---
$ws.Range("K2").Value = -9.508963367108338
$ws.Range("K3").Value = -20.23832653325562
$ws.Range("K4").Value = -16.26923610148342
$ws.Range("K5").Value = -38.87801696020875
